# Apply cell value updates to sheet1 (the active/only worksheet) as described
# by the source diff. All target cells are numeric odds/score values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("J5").Value = 1.11
$ws.Range("K5").Value = 6.5
$ws.Range("W5").Value = 26
$ws.Range("Z5").Value = 6.5

# Row 6
$ws.Range("G6").Value = 1.5
$ws.Range("H6").Value = 4.33
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 1.05
$ws.Range("K6").Value = 11
$ws.Range("L6").Value = 1.33
$ws.Range("M6").Value = 3.25
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 1.8
$ws.Range("P6").Value = 1.4
$ws.Range("Q6").Value = 2.75
$ws.Range("R6").Value = 2.1
$ws.Range("S6").Value = 1.67
$ws.Range("T6").Value = 6
$ws.Range("U6").Value = 6.5
$ws.Range("V6").Value = 8.5
$ws.Range("W6").Value = 10
$ws.Range("X6").Value = 13
$ws.Range("Y6").Value = 29
$ws.Range("Z6").Value = 11
$ws.Range("AA6").Value = 8.5
$ws.Range("AB6").Value = 23
$ws.Range("AC6").Value = 81
$ws.Range("AD6").Value = 12
$ws.Range("AE6").Value = 29
$ws.Range("AF6").Value = 17
$ws.Range("AG6").Value = 67
$ws.Range("AH6").Value = 41
$ws.Range("AI6").Value = 51
$ws.Range("AJ6").Value = 501

# Row 7
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 3.5
$ws.Range("I7").Value = 1.8
$ws.Range("J7").Value = 1.06
$ws.Range("K7").Value = 10
$ws.Range("L7").Value = 1.3
$ws.Range("M7").Value = 3.4
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 1.8
$ws.Range("P7").Value = 1.4
$ws.Range("Q7").Value = 2.75
$ws.Range("S7").Value = 1.83
$ws.Range("T7").Value = 11
$ws.Range("U7").Value = 21
$ws.Range("V7").Value = 13
$ws.Range("W7").Value = 41
$ws.Range("X7").Value = 34
$ws.Range("Y7").Value = 41
$ws.Range("Z7").Value = 10
$ws.Range("AA7").Value = 7
$ws.Range("AC7").Value = 51
$ws.Range("AD7").Value = 7
$ws.Range("AE7").Value = 8.5
$ws.Range("AF7").Value = 8.5
$ws.Range("AJ7").Value = 301

# Row 13
$ws.Range("G13").Value = 3.55
$ws.Range("H13").Value = 2.8
$ws.Range("I13").Value = 2.27
$ws.Range("K13").Value = 6.1
$ws.Range("M13").Value = 2.85
$ws.Range("P13").Value = 1.42
$ws.Range("Q13").Value = 2.67
$ws.Range("S13").Value = 1.95
$ws.Range("Z13").Value = 6.1
$ws.Range("AA13").Value = 5.4
$ws.Range("AD13").Value = 6.9
$ws.Range("AF13").Value = 8.75
$ws.Range("AH13").Value = 19.5

# Row 14
$ws.Range("G14").Value = 2.2
$ws.Range("I14").Value = 2.8
$ws.Range("O14").Value = 2.1
$ws.Range("T14").Value = 8.75
$ws.Range("U14").Value = 10.75
$ws.Range("V14").Value = 7.6
$ws.Range("W14").Value = 18.5
$ws.Range("X14").Value = 13
$ws.Range("Y14").Value = 16.5
$ws.Range("AA14").Value = 6.1
$ws.Range("AH14").Value = 17
$ws.Range("AI14").Value = 19.5

# Row 16
$ws.Range("G16").Value = 1.98
$ws.Range("H16").Value = 3.1
$ws.Range("I16").Value = 3.55
$ws.Range("T16").Value = 5.5
$ws.Range("U16").Value = 7.5
$ws.Range("V16").Value = 7.3
$ws.Range("W16").Value = 14
$ws.Range("X16").Value = 14
$ws.Range("AA16").Value = 5.3
$ws.Range("AB16").Value = 12
$ws.Range("AC16").Value = 55
$ws.Range("AD16").Value = 8.25
$ws.Range("AE16").Value = 15.5
$ws.Range("AF16").Value = 10.25
$ws.Range("AG16").Value = 40
$ws.Range("AH16").Value = 26
$ws.Range("AI16").Value = 32
$ws.Range("AJ16").Value = 400

# Row 17
$ws.Range("G17").Value = 2.07
$ws.Range("H17").Value = 3.2
$ws.Range("I17").Value = 3.55
$ws.Range("J17").Value = 1.06
$ws.Range("K17").Value = 7.1
$ws.Range("L17").Value = 1.29
$ws.Range("M17").Value = 3.25
$ws.Range("N17").Value = 1.87
$ws.Range("O17").Value = 1.83
$ws.Range("P17").Value = 1.42
$ws.Range("Q17").Value = 2.65
$ws.Range("R17").Value = 1.7
$ws.Range("S17").Value = 2.05
$ws.Range("T17").Value = 8
$ws.Range("U17").Value = 10.75
$ws.Range("V17").Value = 8.25
$ws.Range("W17").Value = 20
$ws.Range("X17").Value = 15.5
$ws.Range("Y17").Value = 24
$ws.Range("Z17").Value = 7.1
$ws.Range("AA17").Value = 6.2
$ws.Range("AB17").Value = 13
$ws.Range("AC17").Value = 55
$ws.Range("AD17").Value = 10.25
$ws.Range("AE17").Value = 19
$ws.Range("AF17").Value = 11.75
$ws.Range("AG17").Value = 50
$ws.Range("AH17").Value = 32
$ws.Range("AI17").Value = 37
$ws.Range("AJ17").Value = 400

# Row 18
$ws.Range("G18").Value = 2.65
$ws.Range("H18").Value = 3
$ws.Range("I18").Value = 2.72
$ws.Range("J18").Value = 1.1
$ws.Range("K18").Value = 5.8
$ws.Range("L18").Value = 1.45
$ws.Range("M18").Value = 2.55
$ws.Range("N18").Value = 2.32
$ws.Range("O18").Value = 1.53
$ws.Range("P18").Value = 1.55
$ws.Range("Q18").Value = 2.3
$ws.Range("R18").Value = 1.98
$ws.Range("S18").Value = 1.75
$ws.Range("T18").Value = 7
$ws.Range("U18").Value = 12
$ws.Range("V18").Value = 10.25
$ws.Range("W18").Value = 30
$ws.Range("X18").Value = 25
$ws.Range("Y18").Value = 40
$ws.Range("Z18").Value = 5.8
$ws.Range("AA18").Value = 5.8
$ws.Range("AB18").Value = 16.5
$ws.Range("AC18").Value = 100
$ws.Range("AD18").Value = 6.9
$ws.Range("AE18").Value = 12.5
$ws.Range("AF18").Value = 10.5
$ws.Range("AG18").Value = 32
$ws.Range("AH18").Value = 27
$ws.Range("AI18").Value = 45
$ws.Range("AJ18").Value = 900

# Row 19
$ws.Range("R19").Value = 1.75

# Row 24
$ws.Range("G24").Value = 2.9
$ws.Range("I24").Value = 2.45
$ws.Range("N24").Value = 2.15
$ws.Range("O24").Value = 1.67
$ws.Range("U24").Value = 13
$ws.Range("AF24").Value = 10

# Row 30
$ws.Range("P30").Value = 1.47

# Row 31
$ws.Range("R31").Value = 1.83
$ws.Range("S31").Value = 1.83

# Row 32
$ws.Range("J32").Value = 1.07
$ws.Range("K32").Value = 9
$ws.Range("N32").Value = 2.25
$ws.Range("O32").Value = 1.62
$ws.Range("R32").Value = 1.91
$ws.Range("S32").Value = 1.8

# Row 34
$ws.Range("G34").Value = 2.55
$ws.Range("I34").Value = 2.6

# Row 35
$ws.Range("N35").Value = 1.85
$ws.Range("O35").Value = 1.95
$ws.Range("R35").Value = 1.67

# Row 36
$ws.Range("G36").Value = 1.67
$ws.Range("H36").Value = 3.8
$ws.Range("I36").Value = 4.75
$ws.Range("J36").Value = 1.04
$ws.Range("K36").Value = 13
$ws.Range("R36").Value = 1.67
$ws.Range("S36").Value = 2.1
$ws.Range("U36").Value = 9
$ws.Range("X36").Value = 13
$ws.Range("AB36").Value = 13
$ws.Range("AD36").Value = 15
$ws.Range("AI36").Value = 34

# Row 37
$ws.Range("R37").Value = 1.67

# Row 38
$ws.Range("R38").Value = 1.75
$ws.Range("T38").Value = 7.5
$ws.Range("U38").Value = 9

# Row 40
$ws.Range("K40").Value = 13

# Row 41
$ws.Range("G41").Value = 1.8
$ws.Range("I41").Value = 4.2
$ws.Range("L41").Value = 1.18
$ws.Range("M41").Value = 4.5
$ws.Range("N41").Value = 1.65
$ws.Range("R41").Value = 1.62
$ws.Range("S41").Value = 2.2
$ws.Range("T41").Value = 9
$ws.Range("U41").Value = 10
$ws.Range("Z41").Value = 15
$ws.Range("AF41").Value = 13
$ws.Range("AG41").Value = 41
$ws.Range("AH41").Value = 29

# Row 42
$ws.Range("P42").Value = 1.33
$ws.Range("R42").Value = 1.75

# Row 43
$ws.Range("I43").Value = 4.33
$ws.Range("R43").Value = 1.7
$ws.Range("Z43").Value = 13
$ws.Range("AB43").Value = 13
$ws.Range("AD43").Value = 15
$ws.Range("AI43").Value = 34
$ws.Range("AJ43").Value = 151

# Row 44
$ws.Range("G44").Value = 1.95
$ws.Range("N44").Value = 1.48
$ws.Range("P44").Value = 1.22
$ws.Range("R44").Value = 1.44
$ws.Range("S44").Value = 2.63

# Row 45
$ws.Range("P45").Value = 1.22

# Row 46
$ws.Range("G46").Value = 2.25
$ws.Range("N46").Value = 1.88
$ws.Range("O46").Value = 1.93
$ws.Range("R46").Value = 1.67

